$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.9458723333333334
$ws.Range("H2").Value = 2.837617
$ws.Range("I2").Value = 0.1873686327665471
$ws.Range("J2").Value = 0.1873686327665471
$ws.Range("M2").Value = 8.430598666666667
$ws.Range("N2").Value = 25.291796
$ws.Range("O2").Value = 0.173137200317126
$ws.Range("P2").Value = 0.1731372003171259
$ws.Range("Q2").Value = 7.974270032236889
$ws.Range("R2").Value = 71.768430290132
$ws.Range("S2").Value = 0.03244048050444767
$ws.Range("T2").Value = 0.03244048050444766
$ws.Range("G3").Value = 0.9458723333333334
$ws.Range("H3").Value = 2.837617
$ws.Range("I3").Value = 0.1873686327665471
$ws.Range("J3").Value = 0.1873686327665471
$ws.Range("O3").Value = 0.4685067724286191
$ws.Range("P3").Value = 0.468506772428619
$ws.Range("Q3").Value = 21.57825994895689
$ws.Range("R3").Value = 194.204339540612
$ws.Range("S3").Value = 0.08778347339181818
$ws.Range("T3").Value = 0.08778347339181815
$ws.Range("G4").Value = 0.9458723333333334
$ws.Range("H4").Value = 2.837617
$ws.Range("I4").Value = 0.1873686327665471
$ws.Range("J4").Value = 0.1873686327665471
$ws.Range("M4").Value = 5.125375333333333
$ws.Range("N4").Value = 15.376126
$ws.Range("O4").Value = 0.1052586145864599
$ws.Range("P4").Value = 0.1052586145864599
$ws.Range("Q4").Value = 4.847950725749111
$ws.Range("R4").Value = 43.631556531742
$ws.Range("S4").Value = 0.01972216270196592
$ws.Range("T4").Value = 0.01972216270196592
$ws.Range("G5").Value = 0.9458723333333334
$ws.Range("H5").Value = 2.837617
$ws.Range("I5").Value = 0.1873686327665471
$ws.Range("J5").Value = 0.1873686327665471
$ws.Range("M5").Value = 12.32411466666667
$ws.Range("N5").Value = 36.972344
$ws.Range("O5").Value = 0.2530974126677951
$ws.Range("P5").Value = 0.2530974126677951
$ws.Range("Q5").Value = 11.65703909602756
$ws.Range("R5").Value = 104.913351864248
$ws.Range("S5").Value = 0.04742251616831532
$ws.Range("T5").Value = 0.04742251616831532
$ws.Range("I6").Value = 0.3891165466060174
$ws.Range("J6").Value = 0.3891165466060174
$ws.Range("M6").Value = 8.430598666666667
$ws.Range("N6").Value = 25.291796
$ws.Range("O6").Value = 0.173137200317126
$ws.Range("P6").Value = 0.1731372003171259
$ws.Range("Q6").Value = 16.560511601288
$ws.Range("R6").Value = 149.044604411592
$ws.Range("S6").Value = 0.06737054947643431
$ws.Range("T6").Value = 0.0673705494764343
$ws.Range("I7").Value = 0.3891165466060174
$ws.Range("J7").Value = 0.3891165466060174
$ws.Range("O7").Value = 0.4685067724286191
$ws.Range("P7").Value = 0.468506772428619
$ws.Range("S7").Value = 0.1823037373489556
$ws.Range("T7").Value = 0.1823037373489555
$ws.Range("I8").Value = 0.3891165466060174
$ws.Range("J8").Value = 0.3891165466060174
$ws.Range("M8").Value = 5.125375333333333
$ws.Range("N8").Value = 15.376126
$ws.Range("O8").Value = 0.1052586145864599
$ws.Range("P8").Value = 0.1052586145864599
$ws.Range("Q8").Value = 10.067949030028
$ws.Range("R8").Value = 90.61154127025199
$ws.Range("S8").Value = 0.04095786860841705
$ws.Range("T8").Value = 0.04095786860841705
$ws.Range("I9").Value = 0.3891165466060174
$ws.Range("J9").Value = 0.3891165466060174
$ws.Range("M9").Value = 12.32411466666667
$ws.Range("N9").Value = 36.972344
$ws.Range("O9").Value = 0.2530974126677951
$ws.Range("P9").Value = 0.2530974126677951
$ws.Range("Q9").Value = 24.208677459632
$ws.Range("R9").Value = 217.878097136688
$ws.Range("S9").Value = 0.09848439117221051
$ws.Range("T9").Value = 0.09848439117221051
$ws.Range("G10").Value = 1.499502333333333
$ws.Range("H10").Value = 4.498507
$ws.Range("I10").Value = 0.2970376573303378
$ws.Range("J10").Value = 0.2970376573303378
$ws.Range("M10").Value = 8.430598666666667
$ws.Range("N10").Value = 25.291796
$ws.Range("O10").Value = 0.173137200317126
$ws.Range("P10").Value = 0.1731372003171259
$ws.Range("Q10").Value = 12.64170237206356
$ws.Range("R10").Value = 113.775321348572
$ws.Range("S10").Value = 0.05142826837893252
$ws.Range("T10").Value = 0.05142826837893251
$ws.Range("G11").Value = 1.499502333333333
$ws.Range("H11").Value = 4.498507
$ws.Range("I11").Value = 0.2970376573303378
$ws.Range("J11").Value = 0.2970376573303378
$ws.Range("O11").Value = 0.4685067724286191
$ws.Range("P11").Value = 0.468506772428619
$ws.Range("Q11").Value = 34.20826469118355
$ws.Range("R11").Value = 307.874382220652
$ws.Range("S11").Value = 0.1391641541255947
$ws.Range("T11").Value = 0.1391641541255947
$ws.Range("G12").Value = 1.499502333333333
$ws.Range("H12").Value = 4.498507
$ws.Range("I12").Value = 0.2970376573303378
$ws.Range("J12").Value = 0.2970376573303378
$ws.Range("M12").Value = 5.125375333333333
$ws.Range("N12").Value = 15.376126
$ws.Range("O12").Value = 0.1052586145864599
$ws.Range("P12").Value = 0.1052586145864599
$ws.Range("Q12").Value = 7.685512271542445
$ws.Range("R12").Value = 69.16961044388199
$ws.Range("S12").Value = 0.03126577229059898
$ws.Range("T12").Value = 0.03126577229059898
$ws.Range("G13").Value = 1.499502333333333
$ws.Range("H13").Value = 4.498507
$ws.Range("I13").Value = 0.2970376573303378
$ws.Range("J13").Value = 0.2970376573303378
$ws.Range("M13").Value = 12.32411466666667
$ws.Range("N13").Value = 36.972344
$ws.Range("O13").Value = 0.2530974126677951
$ws.Range("P13").Value = 0.2530974126677951
$ws.Range("Q13").Value = 18.48003869893422
$ws.Range("R13").Value = 166.320348290408
$ws.Range("S13").Value = 0.07517946253521163
$ws.Range("T13").Value = 0.07517946253521163
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.6384806666666667
$ws.Range("H14").Value = 1.915442
$ws.Range("I14").Value = 0.1264771632970977
$ws.Range("J14").Value = 0.1264771632970977
$ws.Range("M14").Value = 8.430598666666667
$ws.Range("N14").Value = 25.291796
$ws.Range("O14").Value = 0.173137200317126
$ws.Range("P14").Value = 0.1731372003171259
$ws.Range("Q14").Value = 5.382774257092445
$ws.Range("R14").Value = 48.444968313832
$ws.Range("S14").Value = 0.02189790195731145
$ws.Range("T14").Value = 0.02189790195731145
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.6384806666666667
$ws.Range("H15").Value = 1.915442
$ws.Range("I15").Value = 0.1264771632970977
$ws.Range("J15").Value = 0.1264771632970977
$ws.Range("O15").Value = 0.4685067724286191
$ws.Range("P15").Value = 0.468506772428619
$ws.Range("Q15").Value = 14.56570967581244
$ws.Range("R15").Value = 131.091387082312
$ws.Range("S15").Value = 0.05925540756225064
$ws.Range("T15").Value = 0.05925540756225063
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.6384806666666667
$ws.Range("H16").Value = 1.915442
$ws.Range("I16").Value = 0.1264771632970977
$ws.Range("J16").Value = 0.1264771632970977
$ws.Range("M16").Value = 5.125375333333333
$ws.Range("N16").Value = 15.376126
$ws.Range("O16").Value = 0.1052586145864599
$ws.Range("P16").Value = 0.1052586145864599
$ws.Range("Q16").Value = 3.272453059743556
$ws.Range("R16").Value = 29.452077537692
$ws.Range("S16").Value = 0.01331281098547796
$ws.Range("T16").Value = 0.01331281098547796
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.6384806666666667
$ws.Range("H17").Value = 1.915442
$ws.Range("I17").Value = 0.1264771632970977
$ws.Range("J17").Value = 0.1264771632970977
$ws.Range("M17").Value = 12.32411466666667
$ws.Range("N17").Value = 36.972344
$ws.Range("O17").Value = 0.2530974126677951
$ws.Range("P17").Value = 0.2530974126677951
$ws.Range("Q17").Value = 7.868708948449778
$ws.Range("R17").Value = 70.818380536048
$ws.Range("S17").Value = 0.03201104279205765
$ws.Range("T17").Value = 0.03201104279205765
